# Update "想去人数" (F column) counts that changed between data refreshes
# for both the "展览" and "全部类型" worksheets (they carry the same data).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 111
    4  = 1612
    6  = 1096
    7  = 16
    8  = 11469
    9  = 28
    15 = 12369
    16 = 13047
    21 = 221
    24 = 108
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
